# Updated symbol list on Fri Feb  3 21:43:58 UTC 2023 with GitHub Actions
# Applies updated price (column D) and Volume(1h) (column E) values for the cryptos sheet.
# Values are assigned with a leading apostrophe to force Excel to store them as text
# (matching the original inlineStr/text cell type), then the cell style is reset to
# "Normal" so no stray number-format/style is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'329.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.26%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'41.14"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'1.85%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.698"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-2.42%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08063"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-0.16%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.031"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'3.67%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'8.706"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-0.91%"
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'-1.84%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.945"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.09%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9240"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-2.20%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1272"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-1.48%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.1941"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-2.20%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'8.266"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-7.85%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09414"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'1.59%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.03695"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'5.49%"
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'9.68%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.001300"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-1.25%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.006250"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.81%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'3.383"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.39%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'-2.54%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'0.02%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2652"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'9.93%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04417"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.08%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'-0.03%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004377"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.57%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001242"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'8.74%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D39").Value = "'0.02821"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'16.73%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05460"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'3.03%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007598"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'1.64%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.009959"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'14.41%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'-0.56%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002133"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'1.36%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.01186"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'7.89%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006728"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-2.42%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-0.04%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.002993"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-5.28%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'33.93%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.04%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.04%"
$ws.Range("E51").Style = "Normal"
